$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) values that look like plain decimal numbers need
# to be forced to Text format first, otherwise Excel auto-converts the
# typed string into a number (losing formatting like trailing zeros)
# -- the source data stores every Price cell as text.

$ws.Range("D2").Value = '29.889.26'
$ws.Range("E2").Value = '  +1.01%  '

$ws.Range("D3").Value = '1.624.05'
$ws.Range("E3").Value = '  +1.11%  '

$ws.Range("E4").Value = '  -0.37%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.14'
$ws.Range("E5").Value = '  +0.84%  '

$ws.Range("E7").Value = '  -0.40%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '29.71'
$ws.Range("E8").Value = '  +11.02%  '

$ws.Range("E9").Value = '  +3.08%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0610'
$ws.Range("E10").Value = '  +1.59%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0913'
$ws.Range("E11").Value = '  +0.61%  '

$ws.Range("D12").Value = '1.857.16'
$ws.Range("E12").Value = '  +1.16%  '

$ws.Range("D13").Value = '1.620.44'
$ws.Range("E13").Value = '  -0.47%  '

$ws.Range("E14").Value = '  +5.95%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.90'
$ws.Range("E15").Value = '  +4.87%  '

$ws.Range("D16").Value = '29.929.24'
$ws.Range("E16").Value = '  +1.14%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '8.80'
$ws.Range("E17").Value = '  +15.82%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '64.67'
$ws.Range("E18").Value = '  +1.96%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '244.32'
$ws.Range("E19").Value = '  +1.50%  '

$ws.Range("D20").Value = '0.0₃0704'
$ws.Range("E20").Value = '  +1.77%  '

$ws.Range("E21").Value = '  -0.28%  '

$ws.Range("E22").Value = '  +3.30%  '

$ws.Range("E23").Value = '  +4.45%  '

$ws.Range("E24").Value = '  +2.50%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '157.08'
$ws.Range("E25").Value = '  +1.60%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.68'
$ws.Range("E26").Value = '  +2.62%  '

$ws.Range("E27").Value = '  +1.77%  '

$ws.Range("E28").Value = '  +3.15%  '

$ws.Range("E29").Value = '  -0.36%  '

$ws.Range("E30").Value = '  +3.38%  '

$ws.Range("E31").Value = '  +5.19%  '

$ws.Range("E32").Value = '  +3.67%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.22'
$ws.Range("E33").Value = '  +3.70%  '

$ws.Range("D34").Value = '1.426.46'
$ws.Range("E34").Value = '  +1.38%  '

$ws.Range("E36").Value = '  -0.34%  '

$ws.Range("E37").Value = '  +2.30%  '

$ws.Range("E38").Value = '  -0.66%  '

$ws.Range("E39").Value = '  +3.10%  '

$ws.Range("E40").Value = '  +3.52%  '

$ws.Range("E41").Value = '  +3.32%  '

$ws.Range("E42").Value = '  +0.47%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.834'
$ws.Range("E43").Value = '  +4.67%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '54.04'
$ws.Range("E44").Value = '  -0.07%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '69.28'
$ws.Range("E45").Value = '  +5.01%  '

$ws.Range("E46").Value = '  +18.37%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.995'
$ws.Range("E47").Value = '  -0.44%  '

$ws.Range("E48").Value = '  +2.80%  '

$ws.Range("D49").Value = '1.765.96'
$ws.Range("E49").Value = '  +1.14%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '88.44'
$ws.Range("E50").Value = '  +2.07%  '

$ws.Range("E51").Value = '  +3.02%  '
